# Updates the cryptos price/volume table to the latest scraped values,
# matching the GitHub Actions "Updated cryptos list" automated commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "68.150.65"; AsText = 0 },
    @{ Cell = "E2"; Value = "  -0.83%  "; AsText = 0 },
    @{ Cell = "D3"; Value = "2.642.70"; AsText = 0 },
    @{ Cell = "E3"; Value = "  -0.48%  "; AsText = 0 },
    @{ Cell = "E4"; Value = "  +0.00%  "; AsText = 0 },
    @{ Cell = "D5"; Value = "596.70"; AsText = 1 },
    @{ Cell = "E5"; Value = "  -0.74%  "; AsText = 0 },
    @{ Cell = "D6"; Value = "156.02"; AsText = 1 },
    @{ Cell = "E6"; Value = "  +0.08%  "; AsText = 0 },
    @{ Cell = "E7"; Value = "  +0.00%  "; AsText = 0 },
    @{ Cell = "D8"; Value = "0.542"; AsText = 1 },
    @{ Cell = "E8"; Value = "  -0.98%  "; AsText = 0 },
    @{ Cell = "D9"; Value = "0.141"; AsText = 1 },
    @{ Cell = "E9"; Value = "  +1.27%  "; AsText = 0 },
    @{ Cell = "E10"; Value = "  -1.21%  "; AsText = 0 },
    @{ Cell = "D11"; Value = "5.24"; AsText = 1 },
    @{ Cell = "E11"; Value = "  +0.07%  "; AsText = 0 },
    @{ Cell = "D12"; Value = "0.351"; AsText = 1 },
    @{ Cell = "E12"; Value = "  -0.11%  "; AsText = 0 },
    @{ Cell = "D13"; Value = "27.99"; AsText = 1 },
    @{ Cell = "E13"; Value = "  -0.01%  "; AsText = 0 },
    @{ Cell = "D14"; Value = "0.0000190"; AsText = 1 },
    @{ Cell = "E14"; Value = "  +0.06%  "; AsText = 0 },
    @{ Cell = "D15"; Value = "3.120.77"; AsText = 0 },
    @{ Cell = "E15"; Value = "  -0.37%  "; AsText = 0 },
    @{ Cell = "D16"; Value = "68.095.01"; AsText = 0 },
    @{ Cell = "E16"; Value = "  -0.74%  "; AsText = 0 },
    @{ Cell = "D17"; Value = "2.637.12"; AsText = 0 },
    @{ Cell = "E17"; Value = "  -0.63%  "; AsText = 0 },
    @{ Cell = "D18"; Value = "11.37"; AsText = 1 },
    @{ Cell = "E18"; Value = "  -0.24%  "; AsText = 0 },
    @{ Cell = "D19"; Value = "362.94"; AsText = 1 },
    @{ Cell = "E19"; Value = "  -0.79%  "; AsText = 0 },
    @{ Cell = "D20"; Value = "7.36"; AsText = 1 },
    @{ Cell = "E20"; Value = "  -0.98%  "; AsText = 0 },
    @{ Cell = "D21"; Value = "4.42"; AsText = 1 },
    @{ Cell = "E21"; Value = "  +3.10%  "; AsText = 0 },
    @{ Cell = "D22"; Value = "4.78"; AsText = 1 },
    @{ Cell = "E22"; Value = "  -3.03%  "; AsText = 0 },
    @{ Cell = "E23"; Value = "  -2.08%  "; AsText = 0 },
    @{ Cell = "D24"; Value = "74.88"; AsText = 1 },
    @{ Cell = "E24"; Value = "  +3.02%  "; AsText = 0 },
    @{ Cell = "E25"; Value = "  -0.03%  "; AsText = 0 },
    @{ Cell = "E26"; Value = "  -3.77%  "; AsText = 0 },
    @{ Cell = "B27"; Value = "Binance-PegBSC-USD"; AsText = 0 },
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; AsText = 0 },
    @{ Cell = "D27"; Value = "1.17"; AsText = 1 },
    @{ Cell = "E27"; Value = "  +16.75%  "; AsText = 0 },
    @{ Cell = "B28"; Value = "WrappedeETH"; AsText = 0 },
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; AsText = 0 },
    @{ Cell = "D28"; Value = "2.773.79"; AsText = 0 },
    @{ Cell = "E28"; Value = "  -0.11%  "; AsText = 0 },
    @{ Cell = "B29"; Value = "PEPE"; AsText = 0 },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; AsText = 0 },
    @{ Cell = "D29"; Value = "0.0000104"; AsText = 1 },
    @{ Cell = "E29"; Value = "  -1.91%  "; AsText = 0 },
    @{ Cell = "D30"; Value = "553.19"; AsText = 1 },
    @{ Cell = "E30"; Value = "  -5.75%  "; AsText = 0 },
    @{ Cell = "D31"; Value = "8.02"; AsText = 1 },
    @{ Cell = "E31"; Value = "  +0.13%  "; AsText = 0 },
    @{ Cell = "E32"; Value = "  -1.73%  "; AsText = 0 },
    @{ Cell = "E33"; Value = "  -1.01%  "; AsText = 0 },
    @{ Cell = "E34"; Value = "  +0.02%  "; AsText = 0 },
    @{ Cell = "E35"; Value = "  -2.11%  "; AsText = 0 },
    @{ Cell = "E36"; Value = "  -0.39%  "; AsText = 0 },
    @{ Cell = "D37"; Value = "161.29"; AsText = 1 },
    @{ Cell = "E37"; Value = "  +0.71%  "; AsText = 0 },
    @{ Cell = "D38"; Value = "19.58"; AsText = 1 },
    @{ Cell = "E38"; Value = "  +1.35%  "; AsText = 0 },
    @{ Cell = "D39"; Value = "0.372"; AsText = 1 },
    @{ Cell = "E39"; Value = "  +0.96%  "; AsText = 0 },
    @{ Cell = "E40"; Value = "  -3.85%  "; AsText = 0 },
    @{ Cell = "D41"; Value = "5.32"; AsText = 1 },
    @{ Cell = "E41"; Value = "  -1.88%  "; AsText = 0 },
    @{ Cell = "D42"; Value = "0.0₆0335"; AsText = 0 },
    @{ Cell = "E42"; Value = "  +4.10%  "; AsText = 0 },
    @{ Cell = "D43"; Value = "17.79"; AsText = 1 },
    @{ Cell = "E43"; Value = "  +0.28%  "; AsText = 0 },
    @{ Cell = "D44"; Value = "2.61"; AsText = 1 },
    @{ Cell = "E44"; Value = "  -2.10%  "; AsText = 0 },
    @{ Cell = "E45"; Value = "  +0.03%  "; AsText = 0 },
    @{ Cell = "D46"; Value = "159.06"; AsText = 1 },
    @{ Cell = "E46"; Value = "  +1.56%  "; AsText = 0 },
    @{ Cell = "D47"; Value = "3.72"; AsText = 1 },
    @{ Cell = "E47"; Value = "  -0.40%  "; AsText = 0 },
    @{ Cell = "D48"; Value = "22.11"; AsText = 1 },
    @{ Cell = "E48"; Value = "  -0.21%  "; AsText = 0 },
    @{ Cell = "E49"; Value = "  -1.63%  "; AsText = 0 },
    @{ Cell = "E50"; Value = "  -0.04%  "; AsText = 0 },
    @{ Cell = "E51"; Value = "  -0.56%  "; AsText = 0 }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText -eq 1) {
        # Force text storage so numeric-looking strings (e.g. "596.70")
        # are not coerced into floating point numbers.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
